$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Datatype EmployeeSituation" block appended below the existing
# "Data String EmployeeSituation" table (rows 39-41), mirroring its
# plain boxed-border style (no fill / no bold).

# Header row (merged across B:C), like the other table headers in the sheet.
$ws.Range("B43").Value = "Datatype EmployeeSituation"
$ws.Range("C43").Value = ""
$ws.Range("B43:C43").Merge()

# Field rows: Type | FieldName
$ws.Range("B44").Value = "String"
$ws.Range("C44").Value = "ID"

$ws.Range("B45").Value = "String"
$ws.Range("C45").Value = "EmployeeSituation"

# Apply the same thin boxed border used by the existing data table (rows 39-41).
$rng = $ws.Range("B43:C45")
$rng.Borders.LineStyle = 1
$rng.Borders.Weight = 2
